{"js": "// Update the date line and the 25 division problems in the practice-sheet\n// table. Each new value is written directly into its own paragraph/table\n// cell by document position (not via a global find-and-replace), because a\n// handful of the \"before\" values repeat in the sheet (e.g. \"72\u00f74=\" appears\n// twice) and some \"after\" values collide with \"before\" values elsewhere in\n// the sheet (e.g. one cell becomes \"51\u00f72=\" while a different cell that\n// already reads \"51\u00f72=\" is changed to something else). Addressing each\n// cell/paragraph directly -- instead of searching for old text -- keeps\n// every substitution unambiguous and leaves the run formatting untouched.\n\nconst body = context.document.body;\n\n// 1) Date heading paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].getRange().insertText(\"2025-05-22 Thursday\", \"Replace\");\n\n// 2) Division problems table: 5 \"data\" rows (table row indices 0, 4, 8,\n// 12, 16 -- the rows in between are blank spacer rows), 5 cells each.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rowReplacements = [\n  { row: 0, values: [\"62\u00f78=\", \"68\u00f73=\", \"30\u00f75=\", \"46\u00f76=\", \"58\u00f75=\"] },\n  { row: 4, values: [\"90\u00f74=\", \"51\u00f72=\", \"73\u00f79=\", \"45\u00f73=\", \"28\u00f75=\"] },\n  { row: 8, values: [\"63\u00f78=\", \"56\u00f76=\", \"46\u00f77=\", \"50\u00f79=\", \"45\u00f75=\"] },\n  { row: 12, values: [\"91\u00f73=\", \"83\u00f79=\", \"79\u00f77=\", \"68\u00f76=\", \"24\u00f75=\"] },\n  { row: 16, values: [\"39\u00f74=\", \"36\u00f78=\", \"63\u00f72=\", \"44\u00f77=\", \"85\u00f79=\"] },\n];\n\nfor (const { row, values } of rowReplacements) {\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    cellParagraphs.items[0].getRange().insertText(values[col], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice-sheet\n# table. Replacements are scoped to their exact paragraph/cell (by\n# document position) rather than a global find-and-replace, because a\n# handful of the \"before\" values repeat (e.g. \"72\u00f74=\" appears twice) and\n# some \"after\" values collide with \"before\" values elsewhere in the sheet\n# (e.g. one cell becomes \"51\u00f72=\" while another cell that already reads\n# \"51\u00f72=\" is changed to something else). Addressing each cell/paragraph\n# directly keeps every substitution unambiguous.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2025-05-22 Thursday\"\n\n# 2) Division problems table: 5 \"data\" rows (Word table row indices 1, 5,\n# 9, 13, 17 -- the rows in between are blank spacer rows), 5 cells each.\n$t = $d.Tables.Item(1)\n\n$rowReplacements = @{\n    1  = @(\"62\u00f78=\", \"68\u00f73=\", \"30\u00f75=\", \"46\u00f76=\", \"58\u00f75=\")\n    5  = @(\"90\u00f74=\", \"51\u00f72=\", \"73\u00f79=\", \"45\u00f73=\", \"28\u00f75=\")\n    9  = @(\"63\u00f78=\", \"56\u00f76=\", \"46\u00f77=\", \"50\u00f79=\", \"45\u00f75=\")\n    13 = @(\"91\u00f73=\", \"83\u00f79=\", \"79\u00f77=\", \"68\u00f76=\", \"24\u00f75=\")\n    17 = @(\"39\u00f74=\", \"36\u00f78=\", \"63\u00f72=\", \"44\u00f77=\", \"85\u00f79=\")\n}\n\nforeach ($row in $rowReplacements.Keys) {\n    $values = $rowReplacements[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
